$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.027.39"
$ws.Range("E2").Value = "  +2.22%  "

# Row 3
$ws.Range("D3").Value = "3.129.55"
$ws.Range("E3").Value = "  +0.46%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "587.15"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "146.95"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.58%  "

# Row 8
$ws.Range("D8").Value = "3.119.46"
$ws.Range("E8").Value = "  +0.68%  "

# Row 9
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.162"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +11.70%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.74"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.468"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.24%  "

# Row 13
$ws.Range("E13").Value = "  +4.09%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "37.22"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.65%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.122"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.00%  "

# Row 16
$ws.Range("D16").Value = "3.646.08"
$ws.Range("E16").Value = "  +0.45%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "63.847.48"
$ws.Range("E17").Value = "  +1.87%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "7.17"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.66%  "

# Row 19
$ws.Range("D19").Value = "3.117.54"
$ws.Range("E19").Value = "  +0.28%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "466.83"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.93%  "

# Row 21
$ws.Range("E21").Value = "  +1.84%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.733"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "

# Row 23
$ws.Range("E23").Value = "  +0.02%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "13.22"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.28%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "82.23"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "

# Row 26
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.98"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +8.91%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.49%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.24"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.09%  "

# Row 30
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("E31").Value = "  +0.16%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "27.07"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "

# Row 33
$ws.Range("E33").Value = "  -2.63%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0884"
$ws.Range("E34").Value = "  +10.15%  "

# Row 35
$ws.Range("E35").Value = "  +7.85%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.06"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.84%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.43"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +14.20%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "6.09"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.17%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "50.92"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.59%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "455.68"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +7.34%  "

# Row 41
$ws.Range("E41").Value = "  -1.07%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0373"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "

# Row 43
$ws.Range("D43").Value = "2.888.99"
$ws.Range("E43").Value = "  -2.07%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.279"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

# Row 45
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.18"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.74%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "35.86"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "124.75"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.54%  "

# Row 50
$ws.Range("E50").Value = "  -0.38%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "24.78"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.50%  "
